$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$xlNone = -4142

# --- Sheet 1 (quality_comparison) ---
# Row 1: give C1 a top+bottom border only (no left/right) and D1 a top+right+bottom
# border (no left), matching the borders already defined in the workbook's style table.
$c1 = $ws1.Range("C1")
$c1.Borders.Item(7).LineStyle = $xlNone   # clear left
$c1.Borders.Item(10).LineStyle = $xlNone  # clear right

$d1 = $ws1.Range("D1")
$d1.Borders.Item(7).LineStyle = $xlNone   # clear left (keeps top/right/bottom)

# Row 2: rename "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2 (computational_comparison) ---
$c1b = $ws2.Range("C1")
$c1b.Borders.Item(7).LineStyle = $xlNone
$c1b.Borders.Item(10).LineStyle = $xlNone

$d1b = $ws2.Range("D1")
$d1b.Borders.Item(7).LineStyle = $xlNone

$f1b = $ws2.Range("F1")
$f1b.Borders.Item(7).LineStyle = $xlNone
$f1b.Borders.Item(10).LineStyle = $xlNone

$g1b = $ws2.Range("G1")
$g1b.Borders.Item(7).LineStyle = $xlNone

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
